# Authentication.xlsx: "Created session key and store in the database"
#
# The old sheet stored two sample credential rows (username/password/salt)
# for "kolonia" and "alex". The new version turns the sheet into a clean
# table header only, adding a 4th column "SessionKey" next to
# Username/Password/Salt, bolding the header row, and dropping the old
# sample data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample data rows (rows 2 & 3: "kolonia" / "alex" records) -
# only the header row should remain.
$ws.Rows("2:3").Delete()

# Add the new 4th column header for the session key.
$ws.Range("D1").Value = "SessionKey"

# Make the header row bold (A1:D1 -> Username, Password, Salt, SessionKey).
$ws.Range("A1:D1").Font.Bold = $true

# Portrait page orientation for printing.
$ws.PageSetup.Orientation = 1

# Matches the saved selection/active-cell state recorded in the sheet.
$ws.Range("D7").Select()
